$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.149.56"
$ws.Range("E2").Value = "  -0.76%  "

$ws.Range("D3").Value = "2.072.72"
$ws.Range("E3").Value = "  -1.32%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "253.81"
$ws.Range("E5").Value = "  +1.13%  "

$ws.Range("D6").Value = "0.677"
$ws.Range("E6").Value = "  +1.95%  "

$ws.Range("D7").Value = "62.22"
$ws.Range("E7").Value = "  +21.56%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "0.393"
$ws.Range("E9").Value = "  +4.77%  "

$ws.Range("D10").Value = "61.64"
$ws.Range("E10").Value = "  +0.12%  "

$ws.Range("D11").Value = "0.0804"
$ws.Range("E11").Value = "  +7.77%  "

$ws.Range("E12").Value = "  +3.10%  "

$ws.Range("D13").Value = "16.44"
$ws.Range("E13").Value = "  +6.71%  "

$ws.Range("D14").Value = "2.373.96"
$ws.Range("E14").Value = "  -1.30%  "

$ws.Range("D15").Value = "0.826"
$ws.Range("E15").Value = "  -0.83%  "

$ws.Range("D16").Value = "5.56"
$ws.Range("E16").Value = "  +8.61%  "

$ws.Range("D17").Value = "2.069.94"
$ws.Range("E17").Value = "  -1.43%  "

$ws.Range("D18").Value = "37.082.24"
$ws.Range("E18").Value = "  -0.46%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "15.54"
$ws.Range("E19").Value = "  +14.46%  "

$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "74.88"
$ws.Range("E20").Value = "  +3.68%  "

$ws.Range("D21").Value = "0.0₃0928"
$ws.Range("E21").Value = "  +11.55%  "

$ws.Range("D22").Value = "5.50"
$ws.Range("E22").Value = "  +5.56%  "

$ws.Range("D23").Value = "240.84"
$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("E24").Value = "  -0.18%  "

$ws.Range("E25").Value = "  -1.37%  "

$ws.Range("D26").Value = "2.39"
$ws.Range("E26").Value = "  +19.59%  "

$ws.Range("D27").Value = "171.79"
$ws.Range("E27").Value = "  +0.99%  "

$ws.Range("D28").Value = "9.36"
$ws.Range("E28").Value = "  +2.00%  "

$ws.Range("E29").Value = "  -1.20%  "

$ws.Range("E30").Value = "  +3.02%  "

$ws.Range("E31").Value = "  +8.20%  "

$ws.Range("D32").Value = "1.11"
$ws.Range("E32").Value = "  +4.93%  "

$ws.Range("E33").Value = "  +5.28%  "

$ws.Range("D34").Value = "4.47"
$ws.Range("E34").Value = "  +9.26%  "

$ws.Range("D35").Value = "0.0900"
$ws.Range("E35").Value = "  -2.55%  "

$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  -0.88%  "

$ws.Range("D38").Value = "1.77"
$ws.Range("E38").Value = "  -3.26%  "

$ws.Range("D39").Value = "0.115"
$ws.Range("E39").Value = "  +27.49%  "

$ws.Range("E40").Value = "  +4.61%  "

$ws.Range("D41").Value = "18.17"
$ws.Range("E41").Value = "  +1.17%  "

$ws.Range("D42").Value = "0.0229"
$ws.Range("E42").Value = "  +1.96%  "

$ws.Range("E43").Value = "  +0.55%  "

$ws.Range("D44").Value = "4.50"
$ws.Range("E44").Value = "  +26.28%  "

$ws.Range("D45").Value = "99.38"
$ws.Range("E45").Value = "  +0.38%  "

$ws.Range("D46").Value = "2.81"
$ws.Range("E46").Value = "  +2.98%  "

$ws.Range("D47").Value = "4.50"
$ws.Range("E47").Value = "  +11.81%  "

$ws.Range("E48").Value = "  +9.61%  "

$ws.Range("D49").Value = "1.309.07"
$ws.Range("E49").Value = "  -0.94%  "

$ws.Range("D50").Value = "2.96"
$ws.Range("E50").Value = "  -2.37%  "

$ws.Range("D51").Value = "6.95"
$ws.Range("E51").Value = "  -0.42%  "
